# Aggiornamento dati Medolla al 23 agosto 2021
# Append nuove righe (344-357) con i dati giornalieri dal 10 al 23 agosto 2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Riga modello da cui copiare la formattazione (data formattata, bordi, allineamento)
$templateRow = 343

# Ciascuna voce: numero riga, data (seriale Excel), nuovi positivi, somma mobile 7gg,
# somma mobile 7gg per 100mila abitanti
$newRows = @(
    @(344, 44418, 2, 13, 208.2999519307803),
    @(345, 44419, 0, 13, 208.2999519307803),
    @(346, 44420, 0, 12, 192.2768787053357),
    @(347, 44421, 6, 16, 256.3691716071143),
    @(348, 44422, 0, 15, 240.3460983816696),
    @(349, 44423, 0, 11, 176.253805479891),
    @(350, 44424, 0, 8, 128.1845858035571),
    @(351, 44425, 0, 6, 96.13843935266785),
    @(352, 44426, 0, 6, 96.13843935266785),
    @(353, 44427, 0, 6, 96.13843935266785),
    @(354, 44428, 1, 1, 16.02307322544464),
    @(355, 44429, 0, 1, 16.02307322544464),
    @(356, 44430, 2, 3, 48.06921967633392),
    @(357, 44431, 0, 3, 48.06921967633392)
)

foreach ($entry in $newRows) {
    $r = $entry[0]
    $data = $entry[1]
    $nuoviPos = $entry[2]
    $sommaMobile = $entry[3]
    $sommaMobile100k = $entry[4]

    # Copia la formattazione (stile data, bordi, allineamento) dalla riga precedente
    $ws.Range("A$($templateRow):D$($templateRow)").Copy()
    $ws.Range("A$($r):D$($r)").PasteSpecial(-4122)

    # Imposta i valori della nuova riga
    $ws.Range("A$r").Value = $data
    $ws.Range("B$r").Value = $nuoviPos
    $ws.Range("C$r").Value = $sommaMobile
    $ws.Range("D$r").Value = $sommaMobile100k
}

$excel.CutCopyMode = $false

Write-Host "Aggiunte righe 344-357 (dati fino al 23 agosto 2021). Nuova dimensione foglio: $($ws.UsedRange.Address())"
